# Generate Report for Handoff
#
# - Flip the handback/status text from "Handed back: in sync with en-US" to
#   "Ready for handoff" everywhere it appears (Overview!E2/F2 and the
#   per-locale "Status" column on the zh-cn / de-de sheets).
# - Bump the handoff-generation timestamps that go along with the new
#   status.
# - Narrow the now much-shorter status columns to fit the shorter text.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet -------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2017-01-03 05:27:00"

# ---- zh-cn sheet ------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2017-01-03 05:26:49"

# ---- de-de sheet --------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2017-01-03 05:27:00"

# ---- Resize the status columns now that the text is shorter ---------
# 16.25 "characters" is the narrowest input that still lands the
# persisted column width at ~17.17 chars (the engine quantizes
# ColumnWidth to 1/6-character steps), matching the new, much narrower
# status columns.
$overview.Columns.Item(5).ColumnWidth = 16.25
$overview.Columns.Item(6).ColumnWidth = 16.25
$zhcn.Columns.Item(3).ColumnWidth = 16.25
$dede.Columns.Item(3).ColumnWidth = 16.25
